$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.058.65'
$ws.Range('E2').Value = '  +12.20%  '

$ws.Range('D3').Value = '1.618.24'
$ws.Range('E3').Value = '  +10.14%  '

$ws.Range('D4').Value = "'0.9891"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -2.00%  '

$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = "'301.02"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +8.76%  '

$ws.Range('B6').Value = 'USDC'
$ws.Range('C6').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D6').Value = "'0.9803"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.47%  '

$ws.Range('D7').Value = "'0.3656"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.70%  '

$ws.Range('D8').Value = "'0.3423"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +11.64%  '

$ws.Range('D9').Value = "'42.55"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +7.90%  '

$ws.Range('D10').Value = "'1.141"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.34%  '

$ws.Range('D11').Value = "'0.07072"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +6.69%  '

$ws.Range('D12').Value = "'0.9809"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.21%  '

$ws.Range('D13').Value = "'20.00"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +10.80%  '

$ws.Range('D14').Value = "'5.886"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +7.83%  '

$ws.Range('D15').Value = "'6.597"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +7.14%  '

$ws.Range('D16').Value = "'0.00001080"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.89%  '

$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '1.612.77'
$ws.Range('E17').Value = '  +9.86%  '

$ws.Range('B18').Value = 'Dai'
$ws.Range('C18').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D18').Value = "'0.9797"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.19%  '

$ws.Range('D19').Value = "'0.06652"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +11.60%  '

$ws.Range('D20').Value = "'78.41"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +13.86%  '

$ws.Range('D21').Value = "'16.13"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +11.10%  '

$ws.Range('D22').Value = "'6.011"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +9.86%  '

$ws.Range('D23').Value = "'11.68"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.53%  '

$ws.Range('D24').Value = '23.047.58'
$ws.Range('E24').Value = '  +12.13%  '

$ws.Range('D25').Value = "'2.370"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.83%  '

$ws.Range('D26').Value = "'2.609"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +25.04%  '

$ws.Range('D27').Value = "'150.05"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.47%  '

$ws.Range('D28').Value = "'19.42"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +13.61%  '

$ws.Range('D29').Value = '1.781.22'
$ws.Range('E29').Value = '  +9.29%  '

$ws.Range('D30').Value = "'125.00"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +9.68%  '

$ws.Range('D31').Value = "'4.078"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.76%  '

$ws.Range('D32').Value = "'6.017"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +22.33%  '

$ws.Range('D33').Value = "'0.9718"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +21.95%  '

$ws.Range('D34').Value = "'1.670"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +15.45%  '

$ws.Range('D35').Value = "'0.08199"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.32%  '

$ws.Range('D36').Value = "'11.82"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +14.50%  '

$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').Value = "'5.171"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +10.07%  '

$ws.Range('B38').Value = 'FraxShare'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D38').Value = "'8.580"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +17.89%  '

$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').Value = "'1.248"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.02%  '

$ws.Range('D40').Value = "'0.06095"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.06%  '

$ws.Range('D41').Value = "'0.02216"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +9.30%  '

$ws.Range('E42').Value = '  +8.03%  '

$ws.Range('D43').Value = "'0.9812"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.31%  '

$ws.Range('D44').Value = "'0.5880"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +12.09%  '

$ws.Range('D45').Value = "'3.770"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.32%  '

$ws.Range('D46').Value = "'12.96"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +7.25%  '

$ws.Range('D47').Value = "'0.5746"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +10.95%  '

$ws.Range('D48').Value = "'125.83"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +6.20%  '

$ws.Range('D49').Value = "'1.969"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +9.37%  '

$ws.Range('D50').Value = "'0.06920"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.54%  '

$ws.Range('D51').Value = "'73.56"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +9.89%  '
